# Generate Report for Handoff
# Adds two new localized files (33c515bb-... and 4b7ac82e-...) that are
# "Ready for handoff" into the Overview / zh-cn / de-de sheets, just
# before the existing "ce5e7eee-..." entry.

$wb = $excel.ActiveWorkbook

$mdBaseUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/13976cec5a4e791122650a917cacd9f7561c3822/e2e"

# -----------------------------------------------------------------
# Sheet 1: "Overview"
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Insert two rows above the existing "ce5e7eee" row (currently row 3),
# pushing it (and the ".localization-config" row after it) down.
$ws1.Rows.Item(3).Resize(2).Insert()

$ws1.Range("A3").Value = "33c515bb-cd31-4a6f-8c16-46bb1b552027.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

$ws1.Range("A4").Value = "4b7ac82e-2cb0-405d-b38d-e1e164478f3e.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"

# Rebuild all hyperlinks on this sheet in the correct left-to-right,
# top-to-bottom order so relationship ids come out sequential again.
$ws1.Hyperlinks.Delete()

$ws1.Hyperlinks.Add($ws1.Range("A2"), "$mdBaseUrl/59b4a092-2a79-4e05-b9e5-24b109c3095c.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "59b4a092-2a79-4e05-b9e5-24b109c3095c.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "$mdBaseUrl/33c515bb-cd31-4a6f-8c16-46bb1b552027.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "33c515bb-cd31-4a6f-8c16-46bb1b552027.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "$mdBaseUrl/4b7ac82e-2cb0-405d-b38d-e1e164478f3e.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "4b7ac82e-2cb0-405d-b38d-e1e164478f3e.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "$mdBaseUrl/ce5e7eee-8057-4cc7-accf-010ede020cce.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "ce5e7eee-8057-4cc7-accf-010ede020cce.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/2fabfc2313b93076f4b6f6fb6383a118a18a4bec/.localization-config", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, ".localization-config") | Out-Null

$ws1.Range("A2:A6").Style = "HyperLink"

# -----------------------------------------------------------------
# Sheet 2: "zh-cn"
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows.Item(3).Resize(2).Insert()

$ws2.Range("A3").Value = "33c515bb-cd31-4a6f-8c16-46bb1b552027.md"
$ws2.Range("B3").Value = "Ready for handoff"
$ws2.Range("C3").Value = "33c515bb-cd31-4a6f-8c16-46bb1b552027.251afffd58fe394e79aaf456db23b5deb709cb31.zh-cn.xlf"
$ws2.Range("D3").Value = "2016-03-09 22:37:01"
$ws2.Range("G3").Value = "0001-01-01 00:00:00"
$ws2.Range("H3").Value = "Include"

$ws2.Range("A4").Value = "4b7ac82e-2cb0-405d-b38d-e1e164478f3e.md"
$ws2.Range("B4").Value = "Ready for handoff"
$ws2.Range("C4").Value = "4b7ac82e-2cb0-405d-b38d-e1e164478f3e.d2e079e0fcbe8039da0e2d3c189a77e0baf9df90.zh-cn.xlf"
$ws2.Range("D4").Value = "2016-03-09 22:37:01"
$ws2.Range("G4").Value = "0001-01-01 00:00:00"
$ws2.Range("H4").Value = "Include"

# Row 6 (the ".localization-config" row, shifted from row 4) also needs
# a "Latest Handoff Datetime" placeholder now that it is not the last row.
$ws2.Range("D6").Value = "0001-01-01 00:00:00"

$ws2.Hyperlinks.Delete()

$ws2.Hyperlinks.Add($ws2.Range("A2"), "$mdBaseUrl/59b4a092-2a79-4e05-b9e5-24b109c3095c.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "59b4a092-2a79-4e05-b9e5-24b109c3095c.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9496262f2b3a4ca279be818ef8473e492a9ae75c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/59b4a092-2a79-4e05-b9e5-24b109c3095c.27aca35bd98b45602e5441f4d58114db17a47550.zh-cn.xlf", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "59b4a092-2a79-4e05-b9e5-24b109c3095c.27aca35bd98b45602e5441f4d58114db17a47550.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/82c1715f978b56e04c76b72a3bff084e324e0ac6/e2e/59b4a092-2a79-4e05-b9e5-24b109c3095c.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "59b4a092-2a79-4e05-b9e5-24b109c3095c.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/529555b2c702fff090834cef671742aad7eadc53/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/59b4a092-2a79-4e05-b9e5-24b109c3095c.27aca35bd98b45602e5441f4d58114db17a47550.zh-cn.xlf", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "59b4a092-2a79-4e05-b9e5-24b109c3095c.27aca35bd98b45602e5441f4d58114db17a47550.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A3"), "$mdBaseUrl/33c515bb-cd31-4a6f-8c16-46bb1b552027.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "33c515bb-cd31-4a6f-8c16-46bb1b552027.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/251afffd58fe394e79aaf456db23b5deb709cb31/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/33c515bb-cd31-4a6f-8c16-46bb1b552027.251afffd58fe394e79aaf456db23b5deb709cb31.zh-cn.xlf", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "33c515bb-cd31-4a6f-8c16-46bb1b552027.251afffd58fe394e79aaf456db23b5deb709cb31.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A4"), "$mdBaseUrl/4b7ac82e-2cb0-405d-b38d-e1e164478f3e.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "4b7ac82e-2cb0-405d-b38d-e1e164478f3e.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d2e079e0fcbe8039da0e2d3c189a77e0baf9df90/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4b7ac82e-2cb0-405d-b38d-e1e164478f3e.d2e079e0fcbe8039da0e2d3c189a77e0baf9df90.zh-cn.xlf", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "4b7ac82e-2cb0-405d-b38d-e1e164478f3e.d2e079e0fcbe8039da0e2d3c189a77e0baf9df90.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A5"), "$mdBaseUrl/ce5e7eee-8057-4cc7-accf-010ede020cce.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "ce5e7eee-8057-4cc7-accf-010ede020cce.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/32d4fd0108c216d0dafbf1e955cc99f298a86ba5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ce5e7eee-8057-4cc7-accf-010ede020cce.d053e05d66fea30943b4118fb366b2017f2a7d5d.zh-cn.xlf", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "ce5e7eee-8057-4cc7-accf-010ede020cce.d053e05d66fea30943b4118fb366b2017f2a7d5d.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/2fabfc2313b93076f4b6f6fb6383a118a18a4bec/.localization-config", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, ".localization-config") | Out-Null

$ws2.Range("A2").Style = "HyperLink"
$ws2.Range("C2").Style = "HyperLink"
$ws2.Range("E2").Style = "HyperLink"
$ws2.Range("F2").Style = "HyperLink"
$ws2.Range("A3:A6").Style = "HyperLink"
$ws2.Range("C3:C5").Style = "HyperLink"

# -----------------------------------------------------------------
# Sheet 3: "de-de"
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows.Item(3).Resize(2).Insert()

$ws3.Range("A3").Value = "33c515bb-cd31-4a6f-8c16-46bb1b552027.md"
$ws3.Range("B3").Value = "Ready for handoff"
$ws3.Range("C3").Value = "33c515bb-cd31-4a6f-8c16-46bb1b552027.251afffd58fe394e79aaf456db23b5deb709cb31.de-de.xlf"
$ws3.Range("D3").Value = "2016-03-09 22:37:06"
$ws3.Range("G3").Value = "0001-01-01 00:00:00"
$ws3.Range("H3").Value = "Include"

$ws3.Range("A4").Value = "4b7ac82e-2cb0-405d-b38d-e1e164478f3e.md"
$ws3.Range("B4").Value = "Ready for handoff"
$ws3.Range("C4").Value = "4b7ac82e-2cb0-405d-b38d-e1e164478f3e.d2e079e0fcbe8039da0e2d3c189a77e0baf9df90.de-de.xlf"
$ws3.Range("D4").Value = "2016-03-09 22:37:06"
$ws3.Range("G4").Value = "0001-01-01 00:00:00"
$ws3.Range("H4").Value = "Include"

$ws3.Range("D6").Value = "0001-01-01 00:00:00"

$ws3.Hyperlinks.Delete()

$ws3.Hyperlinks.Add($ws3.Range("A2"), "$mdBaseUrl/59b4a092-2a79-4e05-b9e5-24b109c3095c.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "59b4a092-2a79-4e05-b9e5-24b109c3095c.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e0b0d2933b0936ea6993d204a62b4d8f2384885b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/59b4a092-2a79-4e05-b9e5-24b109c3095c.27aca35bd98b45602e5441f4d58114db17a47550.de-de.xlf", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "59b4a092-2a79-4e05-b9e5-24b109c3095c.27aca35bd98b45602e5441f4d58114db17a47550.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/b328e7fada08109a24f79d230ddd38847fbdf411/e2e/59b4a092-2a79-4e05-b9e5-24b109c3095c.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "59b4a092-2a79-4e05-b9e5-24b109c3095c.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6d291dcc4ca043420d573d6db00d4f82576237ec/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/59b4a092-2a79-4e05-b9e5-24b109c3095c.27aca35bd98b45602e5441f4d58114db17a47550.de-de.xlf", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "59b4a092-2a79-4e05-b9e5-24b109c3095c.27aca35bd98b45602e5441f4d58114db17a47550.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A3"), "$mdBaseUrl/33c515bb-cd31-4a6f-8c16-46bb1b552027.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "33c515bb-cd31-4a6f-8c16-46bb1b552027.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/251afffd58fe394e79aaf456db23b5deb709cb31/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/33c515bb-cd31-4a6f-8c16-46bb1b552027.251afffd58fe394e79aaf456db23b5deb709cb31.de-de.xlf", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "33c515bb-cd31-4a6f-8c16-46bb1b552027.251afffd58fe394e79aaf456db23b5deb709cb31.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A4"), "$mdBaseUrl/4b7ac82e-2cb0-405d-b38d-e1e164478f3e.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "4b7ac82e-2cb0-405d-b38d-e1e164478f3e.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d2e079e0fcbe8039da0e2d3c189a77e0baf9df90/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4b7ac82e-2cb0-405d-b38d-e1e164478f3e.d2e079e0fcbe8039da0e2d3c189a77e0baf9df90.de-de.xlf", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "4b7ac82e-2cb0-405d-b38d-e1e164478f3e.d2e079e0fcbe8039da0e2d3c189a77e0baf9df90.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A5"), "$mdBaseUrl/ce5e7eee-8057-4cc7-accf-010ede020cce.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "ce5e7eee-8057-4cc7-accf-010ede020cce.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/73c3e57da47597ddab3bffb5be9399fbcf31abe0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ce5e7eee-8057-4cc7-accf-010ede020cce.d053e05d66fea30943b4118fb366b2017f2a7d5d.de-de.xlf", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "ce5e7eee-8057-4cc7-accf-010ede020cce.d053e05d66fea30943b4118fb366b2017f2a7d5d.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/2fabfc2313b93076f4b6f6fb6383a118a18a4bec/.localization-config", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, ".localization-config") | Out-Null

$ws3.Range("A2").Style = "HyperLink"
$ws3.Range("C2").Style = "HyperLink"
$ws3.Range("E2").Style = "HyperLink"
$ws3.Range("F2").Style = "HyperLink"
$ws3.Range("A3:A6").Style = "HyperLink"
$ws3.Range("C3:C5").Style = "HyperLink"
